$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (B1:F1 renamed, G1 new)
# Copy the formatting of the existing header cell (F1) onto the new G1
# header cell before setting its text, so it matches the bold/bordered/
# centered style used by the rest of row 1.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1").PasteSpecial(-4122) | Out-Null

$ws.Range("C1").Value = "AI-Synonyms"
$ws.Range("D1").Value = "Product-AI"
$ws.Range("E1").Value = "Business-Process-AI"
$ws.Range("F1").Value = "Data"
$ws.Range("G1").Value = "Adjectives"

# Data rows: year labels stay the same (2010-2020), values updated and new
# column G (Adjectives) populated for each row.
$data = @(
    @{ Row = 2;  B = 52; C = 2;  D = 1; E = 0; F = 5;  G = 5  },
    @{ Row = 3;  B = 44; C = 5;  D = 4; E = 3; F = 8;  G = 3  },
    @{ Row = 4;  B = 52; C = 1;  D = 1; E = 1; F = 9;  G = 5  },
    @{ Row = 5;  B = 76; C = 3;  D = 1; E = 2; F = 21; G = 4  },
    @{ Row = 6;  B = 72; C = 5;  D = 4; E = 2; F = 15; G = 6  },
    @{ Row = 7;  B = 70; C = 3;  D = 3; E = 2; F = 32; G = 9  },
    @{ Row = 8;  B = 76; C = 1;  D = 3; E = 2; F = 38; G = 18 },
    @{ Row = 9;  B = 81; C = 5;  D = 3; E = 1; F = 46; G = 13 },
    @{ Row = 10; B = 77; C = 8;  D = 4; E = 2; F = 36; G = 23 },
    @{ Row = 11; B = 57; C = 8;  D = 0; E = 1; F = 30; G = 29 },
    @{ Row = 12; B = 64; C = 17; D = 2; E = 1; F = 41; G = 44 }
)

foreach ($r in $data) {
    $row = $r.Row
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}
